$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$lo = $ws.ListObjects.Item(1)

# Extend the table ("表1") by four rows so ref/autoFilter grow from
# A1:R47 to A1:R51 (mirrors adding new rows through the table UI).
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Duplicate the last 4 parameter rows (44:47, the wl=5 CSI500/CSI300 x
# run=1 combos) down into the newly added rows 48:51 - carries over
# values + number formats.
$ws.Range("A44:R47").Copy($ws.Range("A48:R51"))

# Rows 44:47 become the "run=0" variant of those combos.
$ws.Range("A44:A47").Value = "0"

# Rows 48:49 (new CSI500 / CSI300, run=0) keep the wl=5 setup but use
# G=5.0 instead of the old 0.001/0.0001.
$ws.Range("A48:A49").Value = "0"
$ws.Range("L48:L49").Value = "5.0"

# Rows 50:51 (new CSI500 / CSI300, run=1) keep G≈2.5 level but tighten
# wei_tole to 1e-3.
$ws.Range("L50:L51").Value = "2.5"
$ws.Range("O50:O51").Value = "1e-3"

# Match the author's final cursor position/selection.
$ws.Range("A50").Select()
